$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('J2').Value = 3385
$ws.Range('J3').Value = 3549
$ws.Range('H4').Value = 1694
$ws.Range('J4').Value = 782
$ws.Range('J5').Value = 279
$ws.Range('J6').Value = 4144
$ws.Range('H7').Value = 26004
$ws.Range('J7').Value = 12139

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('J2').Value = 48
$ws.Range('J3').Value = 37
$ws.Range('J7').Value = 142

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('J3').Value = 134
$ws.Range('J7').Value = 379

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('J3').Value = 187
$ws.Range('J6').Value = 126
$ws.Range('J7').Value = 444

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('J4').Value = 10
$ws.Range('J7').Value = 92

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('J2').Value = 34
$ws.Range('J3').Value = 30
$ws.Range('J7').Value = 107

$ws = $wb.Worksheets.Item('New City')
$ws.Range('J2').Value = 90
$ws.Range('J7').Value = 315

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('J2').Value = 86
$ws.Range('J7').Value = 369
$ws.Range('J8').Value = 795
$ws.Range('J11').Value = 187
$ws.Range('J12').Value = 25
$ws.Range('J15').Value = 140
$ws.Range('J18').Value = 120
$ws.Range('J20').Value = 253
$ws.Range('J22').Value = 27
$ws.Range('J23').Value = 121
$ws.Range('J24').Value = 37
$ws.Range('J27').Value = 73
$ws.Range('J29').Value = 693
$ws.Range('J31').Value = 92
$ws.Range('J33').Value = 550
$ws.Range('J36').Value = 178
$ws.Range('J37').Value = 379
$ws.Range('J39').Value = 6
$ws.Range('J42').Value = 471
$ws.Range('J48').Value = 121
$ws.Range('J50').Value = 69
$ws.Range('J52').Value = 337
$ws.Range('J54').Value = 231
$ws.Range('J63').Value = 62
$ws.Range('J64').Value = 82
$ws.Range('J65').Value = 315
$ws.Range('J66').Value = 32
$ws.Range('J67').Value = 444
$ws.Range('J69').Value = 30
$ws.Range('J71').Value = 44
$ws.Range('J77').Value = 104
$ws.Range('J79').Value = 361
$ws.Range('J83').Value = 281
$ws.Range('J84').Value = 107
$ws.Range('J85').Value = 551
$ws.Range('J88').Value = 129
$ws.Range('J89').Value = 142
$ws.Range('J90').Value = 142
$ws.Range('J92').Value = 38
$ws.Range('H97').Value = 210
$ws.Range('J97').Value = 76
$ws.Range('J100').Value = 24
$ws.Range('H101').Value = 26004
$ws.Range('J101').Value = 12139

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('J3').Value = 105
$ws.Range('J7').Value = 281

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('J2').Value = 147
$ws.Range('J3').Value = 179
$ws.Range('J5').Value = 23
$ws.Range('J6').Value = 175
$ws.Range('J7').Value = 550

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('J3').Value = 46
$ws.Range('J7').Value = 231

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('J2').Value = 214
$ws.Range('J3').Value = 237
$ws.Range('J6').Value = 175
$ws.Range('J7').Value = 693

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('J3').Value = 23
$ws.Range('J7').Value = 121

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('J2').Value = 137
$ws.Range('J3').Value = 208
$ws.Range('J6').Value = 153
$ws.Range('J7').Value = 551

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('J6').Value = 233
$ws.Range('J7').Value = 471

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range('J3').Value = 12
$ws.Range('J7').Value = 37

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('J3').Value = 42
$ws.Range('J6').Value = 27
$ws.Range('J7').Value = 121

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range('J2').Value = 8
$ws.Range('J7').Value = 30

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('J2').Value = 105
$ws.Range('J3').Value = 134
$ws.Range('J7').Value = 361

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('J2').Value = 25
$ws.Range('J7').Value = 82

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('J6').Value = 62
$ws.Range('J7').Value = 253

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('J6').Value = 64
$ws.Range('J7').Value = 120

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('J2').Value = 64
$ws.Range('J7').Value = 178

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range('J6').Value = 11
$ws.Range('J7').Value = 24

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('J2').Value = 76
$ws.Range('J3').Value = 97
$ws.Range('J7').Value = 337

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('J3').Value = 39
$ws.Range('J7').Value = 140

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('J3').Value = 21
$ws.Range('J7').Value = 69

$ws = $wb.Worksheets.Item('Greektown')
$ws.Range('J5').Value = 3
$ws.Range('J6').Value = 6

$ws = $wb.Worksheets.Item('North Center')
$ws.Range('J6').Value = 17
$ws.Range('J7').Value = 32

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('J2').Value = 68
$ws.Range('J3').Value = 42
$ws.Range('J7').Value = 187

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('J2').Value = 27
$ws.Range('J7').Value = 86

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('H4').Value = 10
$ws.Range('J6').Value = 45
$ws.Range('H7').Value = 210
$ws.Range('J7').Value = 76

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range('J6').Value = 13
$ws.Range('J7').Value = 38

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('J2').Value = 30
$ws.Range('J7').Value = 129

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('J2').Value = 235
$ws.Range('J3').Value = 246
$ws.Range('J6').Value = 247
$ws.Range('J7').Value = 795

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('J6').Value = 27
$ws.Range('J7').Value = 73

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('J6').Value = 44
$ws.Range('J7').Value = 142

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range('J2').Value = 15
$ws.Range('J7').Value = 27

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range('J2').Value = 12
$ws.Range('J7').Value = 44

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('J2').Value = 35
$ws.Range('J3').Value = 36
$ws.Range('J7').Value = 104

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('J6').Value = 122
$ws.Range('J7').Value = 369

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range('J3').Value = 3
$ws.Range('J7').Value = 25
